$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store an exact text string even when the
# text looks like a number (e.g. "335.78"), matching the inline-string
# cell type used in the workbook, without leaving a stray style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "42.996.40"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.387.04"
$ws.Range("E3").Value = "  +4.63%  "

$ws.Range("E4").Value = "  -0.28%  "

Set-TextValue $ws.Range("D5") "335.78"
$ws.Range("E5").Value = "  +8.72%  "

Set-TextValue $ws.Range("D6") "102.01"
$ws.Range("E6").Value = "  -9.50%  "

Set-TextValue $ws.Range("D7") "0.642"
$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("E8").Value = "  +0.07%  "

Set-TextValue $ws.Range("D9") "0.635"
$ws.Range("E9").Value = "  +3.56%  "

Set-TextValue $ws.Range("D10") "41.16"
$ws.Range("E10").Value = "  -7.10%  "

Set-TextValue $ws.Range("D11") "0.0931"
$ws.Range("E11").Value = "  +0.68%  "

Set-TextValue $ws.Range("D12") "8.57"
$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("E13").Value = "  -4.43%  "

Set-TextValue $ws.Range("D14") "16.89"
$ws.Range("E14").Value = "  +8.83%  "

$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").Value = "2.749.91"
$ws.Range("E16").Value = "  +4.89%  "

$ws.Range("D17").Value = "2.397.33"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").Value = "43.015.75"
$ws.Range("E18").Value = "  +0.20%  "

Set-TextValue $ws.Range("D19") "7.52"
$ws.Range("E19").Value = "  +4.69%  "

$ws.Range("E20").Value = "  +0.04%  "

Set-TextValue $ws.Range("D21") "3.88"
$ws.Range("E21").Value = "  +7.66%  "

Set-TextValue $ws.Range("D22") "76.86"
$ws.Range("E22").Value = "  +0.18%  "

Set-TextValue $ws.Range("D23") "272.94"
$ws.Range("E23").Value = "  +6.89%  "

Set-TextValue $ws.Range("D24") "2.38"
$ws.Range("E24").Value = "  -3.86%  "

Set-TextValue $ws.Range("D25") "9.66"
$ws.Range("E25").Value = "  +7.63%  "

Set-TextValue $ws.Range("D26") "11.79"
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  -0.01%  "

Set-TextValue $ws.Range("D28") "24.13"
$ws.Range("E28").Value = "  +8.58%  "

$ws.Range("E29").Value = "  -1.46%  "

Set-TextValue $ws.Range("D30") "174.20"
$ws.Range("E30").Value = "  -0.21%  "

Set-TextValue $ws.Range("D31") "3.12"
$ws.Range("E31").Value = "  -2.03%  "

Set-TextValue $ws.Range("D32") "36.44"
$ws.Range("E32").Value = "  -5.49%  "

Set-TextValue $ws.Range("D33") "0.0920"
$ws.Range("E33").Value = "  +2.15%  "

Set-TextValue $ws.Range("D34") "6.12"
$ws.Range("E34").Value = "  +7.07%  "

Set-TextValue $ws.Range("D35") "0.134"
$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("E36").Value = "  -4.68%  "

Set-TextValue $ws.Range("D37") "4.00"
$ws.Range("E37").Value = "  -4.62%  "

Set-TextValue $ws.Range("D38") "0.0363"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("E39").Value = "  +4.44%  "

Set-TextValue $ws.Range("D40") "2.81"
$ws.Range("E40").Value = "  +10.38%  "

Set-TextValue $ws.Range("D41") "1.53"
$ws.Range("E41").Value = "  +11.10%  "

Set-TextValue $ws.Range("D42") "0.233"
$ws.Range("E42").Value = "  +0.91%  "

Set-TextValue $ws.Range("D43") "69.83"
$ws.Range("E43").Value = "  -4.02%  "

$ws.Range("E44").Value = "  +0.08%  "

Set-TextValue $ws.Range("D45") "92.00"
$ws.Range("E45").Value = "  +44.80%  "

Set-TextValue $ws.Range("D46") "117.58"
$ws.Range("E46").Value = "  +8.33%  "

Set-TextValue $ws.Range("D47") "12.06"
$ws.Range("E47").Value = "  -3.72%  "

Set-TextValue $ws.Range("D48") "5.53"
$ws.Range("E48").Value = "  -2.91%  "

Set-TextValue $ws.Range("D49") "9.09"
$ws.Range("E49").Value = "  +2.51%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.617.23"
$ws.Range("E50").Value = "  +9.00%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D51") "1.29"
$ws.Range("E51").Value = "  -1.11%  "
